# NN needs to be corrected
# Update the computed/predicted statistics for rows 2-9 (columns B-F, H)
# on the active sheet to reflect the corrected neural-network output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1260.862915039062
$ws.Range("C2").Value = 0.9459
$ws.Range("D2").Value = 0.916100025177002
$ws.Range("E2").Value = 1.18560004234314
$ws.Range("F2").Value = 0.8449000120162964
$ws.Range("H2").Value = 0.8559

$ws.Range("B3").Value = 1141.188232421875
$ws.Range("C3").Value = 0.913
$ws.Range("D3").Value = 0.9068000000000001
$ws.Range("E3").Value = 1.23769998550415
$ws.Range("F3").Value = 0.8561000227928162
$ws.Range("H3").Value = 0.7727000000000001

$ws.Range("B4").Value = 764.6358032226562
$ws.Range("C4").Value = 0.8933
$ws.Range("D4").Value = 0.8933
$ws.Range("E4").Value = 0.9495999813079834
$ws.Range("F4").Value = 0.8312000036239624
$ws.Range("H4").Value = 0.6536999999999999

$ws.Range("B5").Value = 848.4271240234375
$ws.Range("C5").Value = 0.895
$ws.Range("D5").Value = 0.8942
$ws.Range("E5").Value = 0.991100013256073
$ws.Range("F5").Value = 0.8418999910354614
$ws.Range("H5").Value = 0.661

$ws.Range("B6").Value = 1145.98583984375
$ws.Range("C6").Value = 0.9066
$ws.Range("D6").Value = 0.9067
$ws.Range("E6").Value = 0.9681000113487244
$ws.Range("F6").Value = 0.8309999704360962
$ws.Range("H6").Value = 0.7723

$ws.Range("B7").Value = 925.951171875
$ws.Range("C7").Value = 0.9325
$ws.Range("D7").Value = 0.9300000071525574
$ws.Range("E7").Value = 1.064100027084351
$ws.Range("F7").Value = 0.8758000135421753
$ws.Range("H7").Value = 0.9787

$ws.Range("B8").Value = 1038.666381835938
$ws.Range("C8").Value = 0.9307
$ws.Range("D8").Value = 0.9292
$ws.Range("E8").Value = 1.075399994850159
$ws.Range("F8").Value = 0.888700008392334
$ws.Range("H8").Value = 0.9718

$ws.Range("B9").Value = 7125.7177734375
$ws.Range("C9").Value = 0.9183
$ws.Range("D9").Value = 0.9121
$ws.Range("E9").Value = 1.23769998550415
$ws.Range("F9").Value = 0.8309999704360962
$ws.Range("H9").Value = 5.6661
